{"js": "// Replace the sole paragraph (\"Te\" + bookmark + \"st document\") with the\n// final four-paragraph structure:\n//   1. \"Test document\" (re-split into \"Test \" / \"document\" runs, with\n//      spell-check proofErr markers around \"document\")\n//   2. an empty paragraph\n//   3. \"Add comments for this documents\" (each word individually wrapped\n//      in spell-check proofErr markers, matching how Word marks words it\n//      does not recognize as the user types them)\n//   4. an empty paragraph that now carries the relocated \"_GoBack\"\n//      bookmark (Word always keeps this bookmark at the most recent\n//      edit position)\n//\n// We build the exact target markup with insertOoxml so the run/proofErr\n// structure matches the authored diff precisely, rather than relying on\n// higher level paragraph/text insertion APIs that wouldn't reproduce the\n// proofErr spell-check tags.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t xml:space=\"preserve\">Test </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>document</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n          </w:p>\n          <w:p/>\n          <w:p>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>Add</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>comments</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>for</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>this</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>documents</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n          </w:p>\n          <w:p>\n            <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n            <w:bookmarkEnd w:id=\"0\"/>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nfirstParagraph.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Replace the whole document body (\"Te\" + bookmark + \"st document\") with\n# the final four-paragraph structure:\n#   1. \"Test document\" (re-split into \"Test \" / \"document\" runs, with\n#      spell-check proofErr markers around \"document\")\n#   2. an empty paragraph\n#   3. \"Add comments for this documents\" (each word individually wrapped\n#      in spell-check proofErr markers, matching how Word marks words it\n#      does not recognize as the user types them)\n#   4. an empty paragraph that now carries the relocated \"_GoBack\"\n#      bookmark (Word always keeps this bookmark at the most recent\n#      edit position)\n#\n# Range.InsertXML replaces the *exact* range's contents with the supplied\n# WordprocessingML, so we target $d.Content (the whole story, including\n# the final paragraph mark) to get full, precise control over the run /\n# proofErr structure -- matching the authored diff exactly, which plain\n# Range.Text / Paragraphs.Add calls could not reproduce (they do not\n# create proofErr spell-check tags).\n\n$d = $word.ActiveDocument\n$r = $d.Content\n\n$ooxml = @\"\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t xml:space=\"preserve\">Test </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>document</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n          </w:p>\n          <w:p/>\n          <w:p>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>Add</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>comments</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>for</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>this</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>documents</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n          </w:p>\n          <w:p>\n            <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n            <w:bookmarkEnd w:id=\"0\"/>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$r.InsertXML($ooxml)\n"}
